$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Despachos")
$ws2 = $wb.Worksheets.Item("NoModificar")

# New dispatch row (row 5) — row 4 stays blank as in the source file.
$ws.Range("A5").Value = "BBV10031256"
$ws.Range("B5").Value = "Sin documento"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "Last Mile"
$ws.Range("G5").Value = "Esteban Gutiérrez"
$ws.Range("H5").Value = 921635782
$ws.Range("J5").Value = "Departamental 1455"
$ws.Range("K5").Value = "La Florida"

# Extend the existing data validations so they keep covering the sheet
# (E/H/I/B/N/M/F) around the new row, and add the two brand-new rules
# (G: non-blank text, J: non-blank text) that ship with the importer.
$ws.Range("E2:E3,E5:E1048576").Validation.Add(1, 1, 3, "0")
$ws.Range("H2:H3,H5:H1048576").Validation.Add(1, 1, 7, "100000000", "999999999")
$ws.Range("I2:I3,I5:I1048576").Validation.Add(7, 1, 3, "=ISNUMBER(MATCH(""*@*.*"",I2,0))")
$ws.Range("G2:G1048576").Validation.Add(6, 1, 5, "1")
$ws.Range("J2:J1048576").Validation.Add(6, 1, 7, "1")

$ws.Range("B2:B3,B5:B1048576").Validation.Add(3, 1, 3, "=NoModificar!$A$1:$A$5")
$ws.Range("N2:N3,N5:N1048576").Validation.Add(3, 1, 3, "=NoModificar!$C$1:$C$2")
$ws.Range("M2:M3,M6:M1048576").Validation.Add(4, 1, 7, "=NoModificar!D1")
$ws.Range("M5").Validation.Add(4, 1, 7, "=NoModificar!D3")
$ws.Range("F2:F3,F5:F1048576").Validation.Add(3, 1, 3, "=NoModificar!$B$1:$B$3")

# Hide the reference sheet now that the importer manages it.
$ws2.Visible = $False

$ws.Range("D15:D16").Select()
